$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "0.530", "1.00")
# are preserved exactly as text, matching the source data (inline strings).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '68.006.48'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '2.521.99'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '595.34'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").Value = '176.24'
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").Value = '2.519.50'
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").Value = '26.61'
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").Value = '2.984.90'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '0.0000178'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '67.860.80'
$ws.Range("E17").Value = '  +1.22%  '
$ws.Range("D18").Value = '2.532.17'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '8.04'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").Value = '366.85'
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '4.67'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '71.33'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = '1.93'
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  +2.20%  '
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = '2.643.58'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '0.0₃0982'
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("D31").Value = '8.33'
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '532.31'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = '1.89'
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").Value = '1.45'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").Value = '157.15'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").Value = '18.82'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").Value = '18.71'
$ws.Range("E40").Value = '  +1.41%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.80'
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").Value = '0.352'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '5.16'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.49'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").Value = '147.18'
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '0.555'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = '3.71'
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0277'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").Value = '0.0754'
$ws.Range("E51").Value = '  -0.84%  '

# Restore default (General) formatting on column D so no extra style is left
# behind from the temporary Text number format applied above.
$dRange.ClearFormats()

Write-Output "cryptos list updated"
